$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-11 Thursday" "2024-04-12 Friday"

Replace-Text "976×2=1952" "575×9=5175"
Replace-Text "275×6=1650" "928×9=8352"
Replace-Text "553×8=4424" "623×6=3738"
Replace-Text "832×8=6656" "503×5=2515"
Replace-Text "899×2=1798" "969×6=5814"
Replace-Text "737×2=1474" "674×6=4044"
Replace-Text "769×4=3076" "158×6=948"
Replace-Text "548×3=1644" "677×4=2708"
Replace-Text "339×8=2712" "391×4=1564"
Replace-Text "438×2=876" "487×4=1948"
Replace-Text "950×8=7600" "815×6=4890"
Replace-Text "278×9=2502" "474×8=3792"
Replace-Text "500×4=2000" "690×9=6210"
Replace-Text "775×9=6975" "719×9=6471"
Replace-Text "448×7=3136" "620×6=3720"
Replace-Text "429×5=2145" "327×6=1962"
Replace-Text "727×2=1454" "431×3=1293"
Replace-Text "231×6=1386" "332×9=2988"
Replace-Text "660×9=5940" "532×4=2128"
Replace-Text "895×3=2685" "915×7=6405"
Replace-Text "814×4=3256" "346×4=1384"
Replace-Text "618×5=3090" "704×2=1408"
Replace-Text "490×8=3920" "434×9=3906"
Replace-Text "730×2=1460" "291×2=582"
Replace-Text "196×8=1568" "557×2=1114"
